$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - edit order: D1, A1, B1, (C1/E1/F1 unchanged), G1
$ws.Range("D1").Value = "Programming"
$ws.Range("A1").Value = "Asp.NetCore"
$ws.Range("B1").Value = "This book is about programming"
$ws.Range("C1").Value = 99
$ws.Range("E1").Value = "Sherwin"
$ws.Range("F1").Value = 1000
$ws.Range("G1").Value = "image1"

# Row 2 - edit order: A2, B2, D2, (C2/F2 unchanged), E2, G2
$ws.Range("A2").Value = "Universe"
$ws.Range("B2").Value = "Universe"
$ws.Range("C2").Value = 199
$ws.Range("D2").Value = "Universe"
$ws.Range("E2").Value = "Paika"
$ws.Range("F2").Value = 2000
$ws.Range("G2").Value = "image2"

# Row 3 (new) - edit order: D3, A3, B3, C3, E3, F3, G3
$ws.Range("D3").Value = "Universe1"
$ws.Range("A3").Value = "Test"
$ws.Range("B3").Value = "Universe"
$ws.Range("C3").Value = 199
$ws.Range("E3").Value = "Paika"
$ws.Range("F3").Value = 2000
$ws.Range("G3").Value = "image3"

$ws.Range("G3").Select()
